$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Camps")

# Insert a column before K - shifts entire column K (and beyond) right by one for ALL rows.
# This preserves K3's pre-existing empty-string cell type by moving it (rather than re-setting it),
# landing it at L3 which is what we want. We'll restore rows 1-2 afterward.
$ws.Range("K3").Insert(-4161)

# Row 1 header: originally K1="Attendees" L1="Camp Committee"; after shift K1 is blank, L1="Attendees", M1="Camp Committee"
# Restore: K1="Attendees", L1="Camp Committee", clear M1
$ws.Range("K1").Value = "Attendees"
$ws.Range("L1").Value = "Camp Committee"
$ws.Range("M1").ClearContents()

# Row 2: originally K2="YCHERN CT113"; after shift K2 blank, L2="YCHERN CT113"
# Target: K2="YCHERN", L2="CT113"
$ws.Range("K2").Value = "YCHERN"
$ws.Range("L2").Value = "CT113"
